# Daily attendance processing - normalize "Recorded By" (column G) ordering
# so that any "System"/"system" entries are listed first, followed by the
# remaining recorder name(s)/email(s) in their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($null -eq $text) { continue }
    if ($text -eq "") { continue }
    if ($text.IndexOf(",") -lt 0) { continue }

    $parts = $text.Split(",")
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p -eq "System" -or $p -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $ordered = $systemParts + $otherParts
    $newText = [string]::Join(", ", $ordered)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
